# Insert a new data row at row 592, shifting all existing rows
# (592..652) down to (593..653), and populate the new row 592 with
# the new weekly price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(592).Insert()

$ws.Cells.Item(592, 1).Value = 3
$ws.Cells.Item(592, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(592, 3).Value = "Coquimbo"
$ws.Cells.Item(592, 4).Value = 44578
$ws.Cells.Item(592, 5).Value = 5
$ws.Cells.Item(592, 6).Value = "Fruta"
$ws.Cells.Item(592, 7).Value = 100102
$ws.Cells.Item(592, 8).Value = "Cítricos"
$ws.Cells.Item(592, 9).Value = 100102004
$ws.Cells.Item(592, 10).Value = "Mandarina"
$ws.Cells.Item(592, 11).Value = "Murcott"
$ws.Cells.Item(592, 12).Value = "Primera"
$ws.Cells.Item(592, 13).Value = 170
$ws.Cells.Item(592, 14).Value = 5500
$ws.Cells.Item(592, 15).Value = 6000
$ws.Cells.Item(592, 16).Value = 5765
$ws.Cells.Item(592, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(592, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(592, 19).Value = 576
$ws.Cells.Item(592, 20).Value = 10
